$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists token-level predictions, grouped by message_id.
# Originally: message_id 0 ("Critical low battery ...", 14 tokens) occupies
# rows 2-15, followed by message_id 1 ("Compass Error ...", 7 tokens) in
# rows 16-22.
# After the edit: message_id 1 ("Compass Error ...") comes first, in rows
# 2-8, followed by message_id 0 ("Critical low battery ...") in rows 9-22.
# All per-token columns (token, token_index, true_label, pred_label, etc.)
# move together with their row, and message_id is 1 throughout.

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'Compass Error Compass Error Compass disconnected .'
$ws.Range("C2").Value = 'Compass'
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 'B-Event'
$ws.Range("F2").Value = 'B-Event'
$ws.Range("G2").Value = $True
$ws.Range("H2").Value = 'Event'
$ws.Range("I2").Value = $True
$ws.Range("J2").Value = 'Event'
$ws.Range("K2").Value = $True
$ws.Range("L2").Value = 'Correct'

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 'Compass Error Compass Error Compass disconnected .'
$ws.Range("C3").Value = 'Error'
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 'E-Event'
$ws.Range("F3").Value = 'E-Event'
$ws.Range("G3").Value = $True
$ws.Range("H3").Value = 'Event'
$ws.Range("I3").Value = $True
$ws.Range("J3").Value = 'Event'
$ws.Range("K3").Value = $True
$ws.Range("L3").Value = 'Correct'

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 'Compass Error Compass Error Compass disconnected .'
$ws.Range("C4").Value = 'Compass'
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 'B-Event'
$ws.Range("F4").Value = 'B-Event'
$ws.Range("G4").Value = $True
$ws.Range("H4").Value = 'Event'
$ws.Range("I4").Value = $True
$ws.Range("J4").Value = 'Event'
$ws.Range("K4").Value = $True
$ws.Range("L4").Value = 'Correct'

# Row 5
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 'Compass Error Compass Error Compass disconnected .'
$ws.Range("C5").Value = 'Error'
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 'E-Event'
$ws.Range("F5").Value = 'E-Event'
$ws.Range("G5").Value = $True
$ws.Range("H5").Value = 'Event'
$ws.Range("I5").Value = $True
$ws.Range("J5").Value = 'Event'
$ws.Range("K5").Value = $True
$ws.Range("L5").Value = 'Correct'

# Row 6
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 'Compass Error Compass Error Compass disconnected .'
$ws.Range("C6").Value = 'Compass'
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 'B-Event'
$ws.Range("F6").Value = 'B-Event'
$ws.Range("G6").Value = $True
$ws.Range("H6").Value = 'Event'
$ws.Range("I6").Value = $True
$ws.Range("J6").Value = 'Event'
$ws.Range("K6").Value = $True
$ws.Range("L6").Value = 'Correct'

# Row 7
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 'Compass Error Compass Error Compass disconnected .'
$ws.Range("C7").Value = 'disconnected'
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 'E-Event'
$ws.Range("F7").Value = 'E-Event'
$ws.Range("G7").Value = $True
$ws.Range("H7").Value = 'Event'
$ws.Range("I7").Value = $True
$ws.Range("J7").Value = 'Event'
$ws.Range("K7").Value = $True
$ws.Range("L7").Value = 'Correct'

# Row 8
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = 'Compass Error Compass Error Compass disconnected .'
$ws.Range("C8").Value = '.'
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 'O'
$ws.Range("F8").Value = 'O'
$ws.Range("G8").Value = $False
$ws.Range("H8").Value = 'None'
$ws.Range("I8").Value = $False
$ws.Range("J8").Value = 'None'
$ws.Range("K8").Value = $True
$ws.Range("L8").Value = 'Correct'

# Row 9
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C9").Value = 'Critical'
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 'B-Event'
$ws.Range("F9").Value = 'B-Event'
$ws.Range("G9").Value = $True
$ws.Range("H9").Value = 'Event'
$ws.Range("I9").Value = $True
$ws.Range("J9").Value = 'Event'
$ws.Range("K9").Value = $True
$ws.Range("L9").Value = 'Correct'

# Row 10
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C10").Value = 'low'
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 'I-Event'
$ws.Range("F10").Value = 'I-Event'
$ws.Range("G10").Value = $True
$ws.Range("H10").Value = 'Event'
$ws.Range("I10").Value = $True
$ws.Range("J10").Value = 'Event'
$ws.Range("K10").Value = $True
$ws.Range("L10").Value = 'Correct'

# Row 11
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C11").Value = 'battery'
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 'E-Event'
$ws.Range("F11").Value = 'E-Event'
$ws.Range("G11").Value = $True
$ws.Range("H11").Value = 'Event'
$ws.Range("I11").Value = $True
$ws.Range("J11").Value = 'Event'
$ws.Range("K11").Value = $True
$ws.Range("L11").Value = 'Correct'

# Row 12
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C12").Value = 'Aircraft'
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 'B-Event'
$ws.Range("F12").Value = 'B-Event'
$ws.Range("G12").Value = $True
$ws.Range("H12").Value = 'Event'
$ws.Range("I12").Value = $True
$ws.Range("J12").Value = 'Event'
$ws.Range("K12").Value = $True
$ws.Range("L12").Value = 'Correct'

# Row 13
$ws.Range("A13").Value = 1
$ws.Range("B13").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C13").Value = 'in'
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 'I-Event'
$ws.Range("F13").Value = 'I-Event'
$ws.Range("G13").Value = $True
$ws.Range("H13").Value = 'Event'
$ws.Range("I13").Value = $True
$ws.Range("J13").Value = 'Event'
$ws.Range("K13").Value = $True
$ws.Range("L13").Value = 'Correct'

# Row 14
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C14").Value = 'Auto'
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 'I-Event'
$ws.Range("F14").Value = 'I-Event'
$ws.Range("G14").Value = $True
$ws.Range("H14").Value = 'Event'
$ws.Range("I14").Value = $True
$ws.Range("J14").Value = 'Event'
$ws.Range("K14").Value = $True
$ws.Range("L14").Value = 'Correct'

# Row 15
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C15").Value = 'Power'
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = 'I-Event'
$ws.Range("F15").Value = 'I-Event'
$ws.Range("G15").Value = $True
$ws.Range("H15").Value = 'Event'
$ws.Range("I15").Value = $True
$ws.Range("J15").Value = 'Event'
$ws.Range("K15").Value = $True
$ws.Range("L15").Value = 'Correct'

# Row 16
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C16").Value = 'Off'
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 'I-Event'
$ws.Range("F16").Value = 'I-Event'
$ws.Range("G16").Value = $True
$ws.Range("H16").Value = 'Event'
$ws.Range("I16").Value = $True
$ws.Range("J16").Value = 'Event'
$ws.Range("K16").Value = $True
$ws.Range("L16").Value = 'Correct'

# Row 17
$ws.Range("A17").Value = 1
$ws.Range("B17").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C17").Value = 'Protection'
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 'E-Event'
$ws.Range("F17").Value = 'E-Event'
$ws.Range("G17").Value = $True
$ws.Range("H17").Value = 'Event'
$ws.Range("I17").Value = $True
$ws.Range("J17").Value = 'Event'
$ws.Range("K17").Value = $True
$ws.Range("L17").Value = 'Correct'

# Row 18
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C18").Value = 'Forced'
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = 'B-Event'
$ws.Range("F18").Value = 'B-Event'
$ws.Range("G18").Value = $True
$ws.Range("H18").Value = 'Event'
$ws.Range("I18").Value = $True
$ws.Range("J18").Value = 'Event'
$ws.Range("K18").Value = $True
$ws.Range("L18").Value = 'Correct'

# Row 19
$ws.Range("A19").Value = 1
$ws.Range("B19").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C19").Value = 'landing'
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 'I-Event'
$ws.Range("F19").Value = 'I-Event'
$ws.Range("G19").Value = $True
$ws.Range("H19").Value = 'Event'
$ws.Range("I19").Value = $True
$ws.Range("J19").Value = 'Event'
$ws.Range("K19").Value = $True
$ws.Range("L19").Value = 'Correct'

# Row 20
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C20").Value = 'in'
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 'I-Event'
$ws.Range("F20").Value = 'I-Event'
$ws.Range("G20").Value = $True
$ws.Range("H20").Value = 'Event'
$ws.Range("I20").Value = $True
$ws.Range("J20").Value = 'Event'
$ws.Range("K20").Value = $True
$ws.Range("L20").Value = 'Correct'

# Row 21
$ws.Range("A21").Value = 1
$ws.Range("B21").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C21").Value = 'progress'
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 'E-Event'
$ws.Range("F21").Value = 'E-Event'
$ws.Range("G21").Value = $True
$ws.Range("H21").Value = 'Event'
$ws.Range("I21").Value = $True
$ws.Range("J21").Value = 'Event'
$ws.Range("K21").Value = $True
$ws.Range("L21").Value = 'Correct'

# Row 22
$ws.Range("A22").Value = 1
$ws.Range("B22").Value = 'Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress .'
$ws.Range("C22").Value = '.'
$ws.Range("D22").Value = 13
$ws.Range("E22").Value = 'O'
$ws.Range("F22").Value = 'O'
$ws.Range("G22").Value = $False
$ws.Range("H22").Value = 'None'
$ws.Range("I22").Value = $False
$ws.Range("J22").Value = 'None'
$ws.Range("K22").Value = $True
$ws.Range("L22").Value = 'Correct'
